$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(817.3, 802.1, 803.4, 803.55, 65, 812.6)
    3  = @(1033.65, 1002, 1032.95, 1029.4, 53, 1003.85)
    4  = @(43897.45, 43655.25, 43800, 43791.85, 14, 43860)
    5  = @(401.4, 394.5, 397.65, 397.8, 96, 401.2)
    6  = @(643.25, 629.15, 636.95, 636.95, 139, 633.65)
    7  = @(511.5, 506.95, 507.7, 507.65, 105, 507.1)
    8  = @(929.1, 924.65, 928.35, 927.8, 124, 927.8)
    9  = @(662.6, 656, 658.15, 657.25, 36, 657.5)
    10 = @(19898, 19815.1, 19840.45, 19841.3, 37, 19827.5)
    11 = @(2390.95, 2367.4, 2381, 2381.5, 75, 2368.05)
    12 = @(567.6, 563.25, 563.65, 563.7, 171, 567.1)
    13 = @(935.15, 928.25, 930.6, 931.15, 11, 929.5)
    14 = @(684.35, 675.9, 682.1, 682.45, 150, 677.8)
    15 = @(126.75, 125.9, 126.65, 126.65, 357, 126.05)
    16 = @(3540.8, 3505.55, 3511, 3511.9, 19, 3537.7)
    17 = @(3412.65, 3368.75, 3403, 3402.25, 20, 3369)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
